$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4147.609
$ws.Range("I64").Value = 3653.4614
$ws.Range("K64").Value = 3653.4614
$ws.Range("M64").Value = -3405.4614

$ws.Range("H67").Value = 4147.609
$ws.Range("I67").Value = 3653.4614
$ws.Range("K67").Value = 3653.4614
$ws.Range("M67").Value = -2795.4614

$ws.Range("H76").Value = 3852
$ws.Range("I76").Value = 3671.9443
$ws.Range("K76").Value = 3671.9443
$ws.Range("M76").Value = -3356.9443

$ws.Range("H79").Value = 3852
$ws.Range("I79").Value = 3671.9443
$ws.Range("K79").Value = 3671.9443
$ws.Range("M79").Value = -2579.9443

$ws.Range("H129").Value = 466.66666
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H137").Value = 3786.0833
$ws.Range("I137").Value = 1949.65
$ws.Range("K137").Value = 5848.950000000001
$ws.Range("M137").Value = -3298.950000000001

$ws.Range("H138").Value = 3187.97
$ws.Range("I138").Value = 1568
$ws.Range("J138").Value = 4060.2615
$ws.Range("K138").Value = 4704
$ws.Range("L138").Value = 12180.7845
$ws.Range("M138").Value = 436
$ws.Range("N138").Value = -22460.7845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5556.276
$ws.Range("I32").Value = 5342.45
$ws.Range("J32").Value = 8000
$ws.Range("K32").Value = 5342.45
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = -5055.45
$ws.Range("N32").Value = -8574

$ws.Range("H44").Value = 1000000000
$ws.Range("J44").Value = 1000000000
$ws.Range("L44").Value = 1000000000
$ws.Range("N44").Value = -1000000976

$ws.Range("H61").Value = 9481.053
$ws.Range("I61").Value = 4706.2085
$ws.Range("J61").Value = 17666.5
$ws.Range("K61").Value = 4706.2085
$ws.Range("L61").Value = 17666.5
$ws.Range("M61").Value = -4494.2085
$ws.Range("N61").Value = -18090.5

$ws.Range("H122").Value = 4311725
$ws.Range("J122").Value = 6945169
$ws.Range("L122").Value = 20835507
$ws.Range("N122").Value = -20840407

$ws.Range("H132").Value = 6352.6763
$ws.Range("I132").Value = 2967.375
$ws.Range("J132").Value = 7394.3076
$ws.Range("K132").Value = 8902.125
$ws.Range("L132").Value = 22182.9228
$ws.Range("M132").Value = -6372.125
$ws.Range("N132").Value = -27242.9228

$ws.Range("H136").Value = 9481.053
$ws.Range("I136").Value = 4706.2085
$ws.Range("J136").Value = 17666.5
$ws.Range("K136").Value = 14118.6255
$ws.Range("L136").Value = 52999.5
$ws.Range("M136").Value = -11568.6255
$ws.Range("N136").Value = -58099.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5335.7095
$ws.Range("I105").Value = 4974.9165
$ws.Range("J105").Value = 6572.7144
$ws.Range("K105").Value = 4974.9165
$ws.Range("L105").Value = 6572.7144
$ws.Range("M105").Value = -3227.9165
$ws.Range("N105").Value = -10066.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 257.85715
$ws.Range("I22").Value = 246.92308
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 246.92308
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 103.07692
$ws.Range("N22").Value = -1100

$ws.Range("H31").Value = 3335.6667
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3335.6667
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3335.6667
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3925.6667

$ws.Range("H34").Value = 3335.6667
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3335.6667
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3335.6667
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3739.6667

$ws.Range("H58").Value = 2460014.8
$ws.Range("I58").Value = 3638609.8
$ws.Range("K58").Value = 3638609.8
$ws.Range("M58").Value = -3638406.8

$ws.Range("H122").Value = 8932.799999999999
$ws.Range("I122").Value = 4811.3335
$ws.Range("J122").Value = 12737.23
$ws.Range("K122").Value = 14434.0005
$ws.Range("L122").Value = 38211.69
$ws.Range("M122").Value = -11984.0005
$ws.Range("N122").Value = -43111.69

$ws.Range("H132").Value = 2877.3447
$ws.Range("I132").Value = 2368.5881
$ws.Range("K132").Value = 7105.7643
$ws.Range("M132").Value = -4575.7643

$ws.Range("H136").Value = 2460014.8
$ws.Range("I136").Value = 3638609.8
$ws.Range("K136").Value = 10915829.4
$ws.Range("M136").Value = -10913279.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 1866.1111
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1866.1111
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 5598.3333
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -7720.3333

$ws.Range("H77").Value = 1866.1111
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1866.1111
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 16794.9999
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -27402.9999

$ws.Range("H80").Value = 3295.6
$ws.Range("I80").Value = 1002
$ws.Range("J80").Value = 3869
$ws.Range("K80").Value = 3006
$ws.Range("L80").Value = 11607
$ws.Range("M80").Value = -2070
$ws.Range("N80").Value = -13479

$ws.Range("H83").Value = 3295.6
$ws.Range("I83").Value = 1002
$ws.Range("J83").Value = 3869
$ws.Range("K83").Value = 9018
$ws.Range("L83").Value = 34821
$ws.Range("M83").Value = -4338
$ws.Range("N83").Value = -44181

$ws.Range("H108").Value = 4138.5
$ws.Range("J108").Value = 7027.5
$ws.Range("L108").Value = 21082.5
$ws.Range("N108").Value = -26842.5

$ws.Range("H123").Value = 6932.857
$ws.Range("I123").Value = 3030
$ws.Range("J123").Value = 7583.3335
$ws.Range("K123").Value = 9090
$ws.Range("L123").Value = 22750.0005
$ws.Range("M123").Value = -6640
$ws.Range("N123").Value = -27650.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 8509.5
$ws.Range("I23").Value = 3012
$ws.Range("J23").Value = 14007
$ws.Range("K23").Value = 3012
$ws.Range("L23").Value = 14007
$ws.Range("M23").Value = -2789
$ws.Range("N23").Value = -14453

$ws.Range("H70").Value = 5470.1606
$ws.Range("I70").Value = 4954
$ws.Range("J70").Value = 5582.3696
$ws.Range("K70").Value = 4954
$ws.Range("L70").Value = 5582.3696
$ws.Range("M70").Value = -4684
$ws.Range("N70").Value = -6122.3696

$ws.Range("H73").Value = 5470.1606
$ws.Range("I73").Value = 4954
$ws.Range("J73").Value = 5582.3696
$ws.Range("K73").Value = 4954
$ws.Range("L73").Value = 5582.3696
$ws.Range("M73").Value = -4018
$ws.Range("N73").Value = -7454.3696

$ws.Range("H122").Value = 9400.714
$ws.Range("I122").Value = 10561
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 31683
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -29233
$ws.Range("N122").Value = -24400

$ws.Range("H132").Value = 43076.184
$ws.Range("I132").Value = 113701.89
$ws.Range("J132").Value = 7763.3335
$ws.Range("K132").Value = 341105.67
$ws.Range("L132").Value = 23290.0005
$ws.Range("M132").Value = -338575.67
$ws.Range("N132").Value = -28350.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 97366
$ws.Range("J106").Value = 97366
$ws.Range("L106").Value = 97366
$ws.Range("N106").Value = -99890

$ws.Range("H122").Value = 6379.9023
$ws.Range("I122").Value = 4987.5386
$ws.Range("J122").Value = 8793.333000000001
$ws.Range("K122").Value = 14962.6158
$ws.Range("L122").Value = 26379.999
$ws.Range("M122").Value = -12512.6158
$ws.Range("N122").Value = -31279.999

$ws.Range("H132").Value = 3755.611
$ws.Range("I132").Value = 3427
$ws.Range("J132").Value = 4272
$ws.Range("K132").Value = 10281
$ws.Range("L132").Value = 12816
$ws.Range("M132").Value = -7751
$ws.Range("N132").Value = -17876

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H74").Value = 8249.888999999999
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 8249.888999999999
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 8249.888999999999
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -10121.889

$ws.Range("H77").Value = 8249.888999999999
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 8249.888999999999
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 24749.667
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -34109.667

$ws.Range("H98").Value = 45000
$ws.Range("J98").Value = 45000
$ws.Range("L98").Value = 45000
$ws.Range("N98").Value = -50990

$ws.Range("H126").Value = 1499.8235
$ws.Range("I126").Value = 1299.7273
$ws.Range("J126").Value = 1866.6666
$ws.Range("K126").Value = 3899.1819
$ws.Range("L126").Value = 5599.9998
$ws.Range("M126").Value = -1429.1819
$ws.Range("N126").Value = -10539.9998

$ws.Range("H132").Value = 2733.1064
$ws.Range("I132").Value = 2130.9143
$ws.Range("J132").Value = 4489.5
$ws.Range("K132").Value = 6392.742899999999
$ws.Range("L132").Value = 13468.5
$ws.Range("M132").Value = -3862.742899999999
$ws.Range("N132").Value = -18528.5

$ws.Range("H136").Value = 4184.9404
$ws.Range("I136").Value = 1833.7778
$ws.Range("J136").Value = 6915.3228
$ws.Range("K136").Value = 5501.3334
$ws.Range("L136").Value = 20745.9684
$ws.Range("M136").Value = -2951.3334
$ws.Range("N136").Value = -25845.9684
